# Quarterly income-statement update for dollar_cumulative.xlsx
#
# The oldest reporting period ("6 ماهه منتهی به 1399/06", published
# 1400-08-30 (4)) is dropped, every remaining period shifts one column to
# the left, and a brand-new period is appended in the now-empty last
# column: "12 ماهه منتهی به 1401/12", published 1402-02-23 (7)/1402-02-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete first data column (D). Excel shifts every column
# to its right (headers, published dates, and all figures) one slot left,
# and shrinks the column-width table accordingly.
$ws.Range("D:D").Delete()

# New last column is now M (13) - give it the same width style as the
# other "wide" (29-char) columns it joins (E, I were already 29-wide).
$ws.Columns.Item(13).ColumnWidth = 28.17

# Clone the full formatting (fill/border/font/alignment) of the column it
# sits next to (L, formerly M) onto the brand-new column M so the new
# period's cells look like every other data column.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)

# New period header + publish date for the newly appended column M.
# Leading "'" forces plain text (matches the source workbook, where every
# publish-date cell - even bare yyyy-mm-dd ones - is stored as text, not
# an actual date) instead of Excel's automatic date auto-recognition.
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value = "'1402-02-23"

# Updated publish-date label for column I (was J before the shift):
# "1401-10-28 (6)" -> "1402-02-23 (7)"
$ws.Cells.Item(9, 9).Value = "1402-02-23 (7)"

# New quarter's figures (column M), row by row.
$ws.Cells.Item(11, 13).Value = 6237391
$ws.Cells.Item(12, 13).Value = -5495390
$ws.Cells.Item(13, 13).Value = 742002
$ws.Cells.Item(14, 13).Value = -39711
$ws.Cells.Item(15, 13).Value = "-"
$ws.Cells.Item(16, 13).Value = 17054
$ws.Cells.Item(17, 13).Value = 719344
$ws.Cells.Item(18, 13).Value = -10894
$ws.Cells.Item(19, 13).Value = 205547
$ws.Cells.Item(20, 13).Value = 913997
$ws.Cells.Item(21, 13).Value = -63785
$ws.Cells.Item(22, 13).Value = 850212
$ws.Cells.Item(23, 13).Value = "-"
$ws.Cells.Item(24, 13).Value = 850212
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 13).Value = 779279
$ws.Cells.Item(27, 13).Value = 0
